# Update "Generate Report for Handback" timestamps

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-10-20 08:43:13"

$wsZhCn.Range("H3").Value = "2016-10-20 08:43:01"
$wsZhCn.Range("K3").Value = "2016-10-20 08:43:46"

$wsDeDe.Range("K3").Value = "2016-10-20 08:44:05"
